# Insert a new record row at row 390 (weekly update adds one new observation),
# pushing the existing rows 390-446 down to 391-447.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(390).Insert()

$ws.Range("A390").Value = 4
$ws.Range("B390").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C390").Value = "Los Lagos"
$ws.Range("D390").Value = 45154
$ws.Range("E390").Value = 10
$ws.Range("F390").Value = 100112032
$ws.Range("G390").Value = "Zapallo italiano"
$ws.Range("H390").Value = "Sin especificar"
$ws.Range("I390").Value = "Primera"
$ws.Range("J390").Value = 70
$ws.Range("K390").Value = 19000
$ws.Range("L390").Value = 19000
$ws.Range("M390").Value = 19000
$ws.Range("N390").Value = "`$/caja 50 unidades"
$ws.Range("O390").Value = "Región de Arica y Parinacota"
$ws.Range("P390").Value = 380
$ws.Range("Q390").Value = 50
$ws.Range("R390").Value = "Hortaliza"
